$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = @'
WITH {
    phs_accession: "phs001713",
    subject_ids: [],
    experimental_strategies: [],
    genders: [],
    sample_tumor_statuses: [],
    file_types: [],
    library_strategies: [],
    library_sources: [],
    library_selections: [],
    library_layouts: [],
    platforms: [],
    instrument_models: [],
    reference_genome_assemblies: [],
    primary_diagnoses: [],
    num_study_samples_min: 0,
    num_study_samples_max: 0,
    num_study_participants_max: 0,
    num_study_participants_min: 0
} AS inputs, "Not specified in data" AS na
MATCH (s:study {phs_accession: inputs.phs_accession})
OPTIONAL MATCH (s)<--(p:participant)
WITH inputs, na, s, count(distinct p) AS num_p
WHERE
    (inputs.num_study_participants_min = 0 OR num_p >= inputs.num_study_participants_min) AND
    (inputs.num_study_participants_max = 0 OR num_p <= inputs.num_study_participants_max)
OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
WITH inputs, na, s, count(distinct samp) AS num_samp
WHERE
    (inputs.num_study_samples_min = 0 OR num_samp >= inputs.num_study_samples_min) AND
    (inputs.num_study_participants_max = 0 OR num_samp <= inputs.num_study_participants_max)
MATCH (s)<--(f:file)
WITH inputs, na, f, {
file_name: coalesce(f.file_name, ""),
file_type: coalesce(f.file_type, ""),
study_name: coalesce(s.study_name, ""),
accession: coalesce(s.phs_accession, "")
} AS data,
apoc.coll.toSet(apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\\s{0,1}")) AS experimental_strategies
WHERE
(size(inputs.experimental_strategies) = 0 OR size(apoc.coll.intersection(inputs.experimental_strategies, experimental_strategies)) > 0) AND
    (size(inputs.file_types) = 0 OR f.file_type IN inputs.file_types)
OPTIONAL MATCH (f)-->(samp)
WITH inputs, na, f, apoc.map.merge(data, {
sample_id: coalesce(apoc.text.join(apoc.coll.sort(collect(distinct samp.sample_id)), ", "), "")
}) AS data,
COLLECT(DISTINCT samp.sample_tumor_status) AS sample_tumor_status
WHERE
(size(inputs.sample_tumor_statuses) = 0 OR size(apoc.coll.intersection(inputs.sample_tumor_statuses, sample_tumor_status)) > 0)
OPTIONAL MATCH (f)-->(:sample)-->(p:participant)
WITH inputs, na, f,
apoc.map.merge(data, {
participant_id: coalesce(apoc.text.join(collect(distinct p.participant_id), ", "), "")
}) AS data,
COLLECT(DISTINCT p.participant_id) AS participant_id,
COLLECT(DISTINCT p.gender) AS gender
WHERE
    (size(inputs.subject_ids) = 0 OR size(apoc.coll.intersection(inputs.participant_id, participant_id)) > 0) AND
    (size(inputs.genders) = 0 OR size(apoc.coll.intersection(inputs.gender, gender)) > 0)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH inputs, na, f, data,
COLLECT(DISTINCT g.library_strategy) AS library_strategy,
COLLECT(DISTINCT g.library_source) AS library_source,
COLLECT(DISTINCT g.library_selection) AS library_selection,
COLLECT(DISTINCT g.library_layout) AS library_layout,
COLLECT(DISTINCT g.platform) AS platform,
COLLECT(DISTINCT g.instrument_model) AS instrument_models,
COLLECT(DISTINCT g.reference_genome_assembly) AS reference_genome_assembly
WHERE
    (size(inputs.library_strategies) = 0 OR size(apoc.coll.intersection(inputs.library_strategies, library_strategy)) > 0) AND
    (size(inputs.library_sources) = 0 OR size(apoc.coll.intersection(inputs.library_sources, library_source)) > 0) AND
    (size(inputs.library_selections) = 0 OR size(apoc.coll.intersection(inputs.library_selections, library_selection)) > 0) AND
    (size(inputs.library_layouts) = 0 OR size(apoc.coll.intersection(inputs.library_layouts, library_layout)) > 0) AND
    (size(inputs.platforms) = 0 OR size(apoc.coll.intersection(inputs.platforms, platform)) > 0) AND
    (size(inputs.instrument_models) = 0 OR size(apoc.coll.intersection(inputs.instrument_models, instrument_models)) > 0) AND
    (size(inputs.reference_genome_assemblies) = 0 OR size(apoc.coll.intersection(inputs.reference_genome_assemblies, reference_genome_assembly)) > 0)
OPTIONAL MATCH (f)-->(:sample)-->(:participant)<--(diag:diagnosis)
WITH inputs, na, f, data,
COLLECT(DISTINCT diag.primary_diagnosis) AS primary_diagnosis
WHERE
    (size(inputs.primary_diagnoses) = 0 OR size(apoc.coll.intersection(inputs.primary_diagnoses, primary_diagnosis)) > 0)
RETURN
data.file_name AS `File Name`,
data.study_name AS `Study Name`,
data.accession AS `Accession`,
data.participant_id AS `Participant Id`,
data.sample_id AS `Sample Id`,
data.file_type AS `File Type`
ORDER BY `File Name`
LIMIT 100
'@
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5
[void]$ws.Range("C3").Select()
